$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matriz_Resultados: the pairwise win/loss matrix is corrected to all
#    zeros (no model is statistically better than another after the DM/HLN
#    correction).
# ---------------------------------------------------------------------------
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("B2:E5").Value = 0

# ---------------------------------------------------------------------------
# 2) P_valores: corrected (symmetric) p-value matrix. Diagonal stays 1.
# ---------------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.01168083775224749
$wsP.Range("D2").Value = 0.02129532017764135
$wsP.Range("E2").Value = 0.04463706537645407

$wsP.Range("B3").Value = 0.01168083775224749
$wsP.Range("D3").Value = 0.02993862536184366
$wsP.Range("E3").Value = 0.05437911888773295

$wsP.Range("B4").Value = 0.02129532017764135
$wsP.Range("C4").Value = 0.02993862536184366
$wsP.Range("E4").Value = 0.07902955602612649

$wsP.Range("B5").Value = 0.04463706537645407
$wsP.Range("C5").Value = 0.05437911888773295
$wsP.Range("D5").Value = 0.07902955602612649

# ---------------------------------------------------------------------------
# 3) Estadisticos_HLN_DM: corrected Harvey-Leybourne-Newbold DM statistics
#    matrix (antisymmetric). Diagonal stays 0.
# ---------------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_HLN_DM")

$wsE.Range("C2").Value = -2.66533750903342
$wsE.Range("D2").Value = -2.414323232006722
$wsE.Range("E2").Value = -2.085129740167291

$wsE.Range("B3").Value = 2.66533750903342
$wsE.Range("D3").Value = -2.265916207543778
$wsE.Range("E3").Value = -1.992628696224204

$wsE.Range("B4").Value = 2.414323232006722
$wsE.Range("C4").Value = 2.265916207543778
$wsE.Range("E4").Value = -1.810714199399733

$wsE.Range("B5").Value = 2.085129740167291
$wsE.Range("C5").Value = 1.992628696224204
$wsE.Range("D5").Value = 1.810714199399733

# ---------------------------------------------------------------------------
# 4) Resumen_Modelos: recomputed summary (Victorias/Derrotas/Empates +
#    Tasa_Victoria now all 0/0/3/"0.0%" since no comparison is significant
#    any more) and the ECRPS text labels re-point at the surviving shared
#    strings.
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Resumen_Modelos")

$wsR.Range("B2:D5").Value = 0
$wsR.Range("D2").Value = 3
$wsR.Range("D3").Value = 3
$wsR.Range("D4").Value = 3
$wsR.Range("D5").Value = 3

# Force the E:H columns to stay plain text (otherwise values such as
# "0.0%" / "1.2610" get auto-coerced into numbers by value-entry type
# inference), then restore the default "Normal" style so no visible
# formatting change is introduced.
$rngText = $wsR.Range("E2:H5")
$rngText.NumberFormat = "@"

$wsR.Range("E2").Value = "0.0%"
$wsR.Range("F2").Value = "1.2610"
$wsR.Range("G2").Value = "1.6257"
$wsR.Range("H2").Value = "1.2892"

$wsR.Range("E3").Value = "0.0%"
$wsR.Range("F3").Value = "1.4231"
$wsR.Range("G3").Value = "2.0524"
$wsR.Range("H3").Value = "1.4422"

$wsR.Range("E4").Value = "0.0%"
$wsR.Range("F4").Value = "2.1049"
$wsR.Range("G4").Value = "4.6894"
$wsR.Range("H4").Value = "2.2279"

$wsR.Range("E5").Value = "0.0%"
$wsR.Range("F5").Value = "3.1269"
$wsR.Range("G5").Value = "11.7628"
$wsR.Range("H5").Value = "3.7617"

$rngText.Style = "Normal"

Write-Host "edit applied"
